$wb = $excel.ActiveWorkbook

# --- Rename headers on existing sheets ---
$ws1 = $wb.Worksheets.Item("Weekly Quantity")
$ws1.Range("B1").Value = "Weekly_PO_Qty"

$ws2 = $wb.Worksheets.Item("Monthly Trend")
$ws2.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet at the end ---
$ws3 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3.Name = "PO Forecast"

# --- Header row ---
$ws3.Range("A1").Value = "ds"
$ws3.Range("B1").Value = "PO_Forecast"
$ws3.Range("C1").Value = "yhat_lower"
$ws3.Range("D1").Value = "yhat_upper"

# Reuse the existing bold header style (s=1) from the Weekly Quantity sheet
$ws1.Range("A1:B1").Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)

# --- Data rows ---
$ws3.Range("A2").Value = 44948.99999999999
$ws3.Range("B2").Value = 11
$ws3.Range("C2").Value = 0.9489199441047816
$ws3.Range("D2").Value = 21.6031900645624
$ws3.Range("A3").Value = 44955.99999999999
$ws3.Range("B3").Value = 12
$ws3.Range("C3").Value = 1.160055549482488
$ws3.Range("D3").Value = 22.80292413757278
$ws3.Range("A4").Value = 45074.99999999999
$ws3.Range("B4").Value = 36
$ws3.Range("C4").Value = 24.77897592688025
$ws3.Range("D4").Value = 45.3973204507783
$ws3.Range("A5").Value = 45081.99999999999
$ws3.Range("B5").Value = 37
$ws3.Range("C5").Value = 26.26166732516216
$ws3.Range("D5").Value = 47.56031376012243
$ws3.Range("A6").Value = 45088.99999999999
$ws3.Range("B6").Value = 38
$ws3.Range("C6").Value = 28.2816225331729
$ws3.Range("D6").Value = 48.88416453672037
$ws3.Range("A7").Value = 45095.99999999999
$ws3.Range("B7").Value = 40
$ws3.Range("C7").Value = 28.94686434399364
$ws3.Range("D7").Value = 50.18398331873819
$ws3.Range("A8").Value = 45102.99999999999
$ws3.Range("B8").Value = 41
$ws3.Range("C8").Value = 31.31853587417132
$ws3.Range("D8").Value = 52.08999007407952
$ws3.Range("A9").Value = 45109.99999999999
$ws3.Range("B9").Value = 43
$ws3.Range("C9").Value = 32.37516515073114
$ws3.Range("D9").Value = 52.81874154053377
$ws3.Range("A10").Value = 45116.99999999999
$ws3.Range("B10").Value = 44
$ws3.Range("C10").Value = 33.62577450380501
$ws3.Range("D10").Value = 54.64268051154951
$ws3.Range("A11").Value = 45123.99999999999
$ws3.Range("B11").Value = 45
$ws3.Range("C11").Value = 34.86419284576928
$ws3.Range("D11").Value = 55.41169060873649
$ws3.Range("A12").Value = 45130.99999999999
$ws3.Range("B12").Value = 47
$ws3.Range("C12").Value = 35.99975895032063
$ws3.Range("D12").Value = 57.02589028015482
$ws3.Range("A13").Value = 45137.99999999999
$ws3.Range("B13").Value = 48
$ws3.Range("C13").Value = 37.48845880330999
$ws3.Range("D13").Value = 58.42830456522274

# Reuse the existing date style (s=2) for column A
$ws1.Range("A2").Copy()
$ws3.Range("A2:A13").PasteSpecial(-4122)

# Restore original active sheet/selection (unchanged by this edit)
$ws1.Activate() | Out-Null
$ws1.Range("A1").Select() | Out-Null
